$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "53.575.15"
$ws.Range("E2").Value = "  +3.97%  "
$ws.Range("D3").Value = "3.148.17"
$ws.Range("E3").Value = "  +2.93%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "397.32"
$ws.Range("E5").Value = "  +2.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.94"
$ws.Range("E6").Value = "  +6.74%  "
$ws.Range("E7").Value = "  +0.78%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +4.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.02"
$ws.Range("E10").Value = "  +5.90%  "
$ws.Range("E11").Value = "  +1.37%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0873"
$ws.Range("E12").Value = "  +1.54%  "
$ws.Range("D13").Value = "3.655.03"
$ws.Range("E13").Value = "  +3.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.11"
$ws.Range("E14").Value = "  +3.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.02"
$ws.Range("E15").Value = "  +3.14%  "
$ws.Range("E16").Value = "  +8.82%  "
$ws.Range("D17").Value = "3.153.67"
$ws.Range("E17").Value = "  +3.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.57"
$ws.Range("E18").Value = "  -0.80%  "
$ws.Range("D19").Value = "53.470.06"
$ws.Range("E19").Value = "  +3.67%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.29"
$ws.Range("E20").Value = "  +4.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.87"
$ws.Range("E21").Value = "  +3.37%  "
$ws.Range("D22").Value = "0.0₃0975"
$ws.Range("E22").Value = "  +0.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.96"
$ws.Range("E23").Value = "  +1.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "271.03"
$ws.Range("E24").Value = "  +1.09%  "
$ws.Range("E25").Value = "  +3.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.10"
$ws.Range("E26").Value = "  -1.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "27.60"
$ws.Range("E27").Value = "  +2.91%  "
$ws.Range("E28").Value = "  +1.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.170"
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.111"
$ws.Range("E31").Value = "  +2.77%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.01"
$ws.Range("E32").Value = "  +7.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "37.22"
$ws.Range("E33").Value = "  +7.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0496"
$ws.Range("E34").Value = "  +10.90%  "
$ws.Range("E35").Value = "  +0.62%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "50.53"
$ws.Range("E36").Value = "  +1.10%  "
$ws.Range("E37").Value = "  +10.13%  "
$ws.Range("E38").Value = "  -0.16%  "
$ws.Range("E39").Value = "  +8.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.14"
$ws.Range("E40").Value = "  +10.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.36"
$ws.Range("E41").Value = "  +2.41%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.291"
$ws.Range("E42").Value = "  -0.69%  "
$ws.Range("E43").Value = "  +1.56%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "130.47"
$ws.Range("E44").Value = "  +4.09%  "
$ws.Range("E45").Value = "  +1.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.23"
$ws.Range("E46").Value = "  +1.42%  "
$ws.Range("E47").Value = "  -1.00%  "
$ws.Range("D48").Value = "2.084.65"
$ws.Range("E48").Value = "  +2.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.39"
$ws.Range("E49").Value = "  -1.67%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0511"
$ws.Range("E50").Value = "  +22.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0338"
$ws.Range("E51").Value = "  +6.07%  "
